$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F ("time_taken"). Copy formatting from the existing header
# cell E1 (bold font, thin border, centered/top alignment) onto F1, then set
# its text so the paste-formats step does not clobber the value.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data rows 2-47: plain (unstyled) time_taken values, written as text so
# Excel does not reinterpret the timestamp-looking strings as dates.
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:41:33.908087"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:41:33.908101"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:41:33.908105"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:41:33.908108"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:41:33.908112"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:41:33.908115"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:41:33.908118"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:41:33.908121"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:41:33.908124"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:41:33.908127"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:41:33.908130"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:41:33.908133"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:41:33.908136"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:41:33.908139"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:41:33.908142"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:41:33.908145"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:41:33.908148"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:41:33.908151"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:41:33.908154"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:41:33.908157"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:41:33.908160"
$ws.Cells.Item(23, 6).Value = "2021-10-05 13:41:33.908163"
$ws.Cells.Item(24, 6).Value = "2021-10-05 13:41:33.908166"
$ws.Cells.Item(25, 6).Value = "2021-10-05 13:41:33.908169"
$ws.Cells.Item(26, 6).Value = "2021-10-05 13:41:33.908172"
$ws.Cells.Item(27, 6).Value = "2021-10-05 13:41:33.908175"
$ws.Cells.Item(28, 6).Value = "2021-10-05 13:41:33.908178"
$ws.Cells.Item(29, 6).Value = "2021-10-05 13:41:33.908181"
$ws.Cells.Item(30, 6).Value = "2021-10-05 13:41:33.908184"
$ws.Cells.Item(31, 6).Value = "2021-10-05 13:41:33.908187"
$ws.Cells.Item(32, 6).Value = "2021-10-05 13:41:33.908190"
$ws.Cells.Item(33, 6).Value = "2021-10-05 13:41:33.908193"
$ws.Cells.Item(34, 6).Value = "2021-10-05 13:41:33.908196"
$ws.Cells.Item(35, 6).Value = "2021-10-05 13:41:33.908199"
$ws.Cells.Item(36, 6).Value = "2021-10-05 13:41:33.908202"
$ws.Cells.Item(37, 6).Value = "2021-10-05 13:41:33.908205"
$ws.Cells.Item(38, 6).Value = "2021-10-05 13:41:33.908208"
$ws.Cells.Item(39, 6).Value = "2021-10-05 13:41:33.908211"
$ws.Cells.Item(40, 6).Value = "2021-10-05 13:41:33.908214"
$ws.Cells.Item(41, 6).Value = "2021-10-05 13:41:33.908217"
$ws.Cells.Item(42, 6).Value = "2021-10-05 13:41:33.908220"
$ws.Cells.Item(43, 6).Value = "2021-10-05 13:41:33.908224"
$ws.Cells.Item(44, 6).Value = "2021-10-05 13:41:33.908226"
$ws.Cells.Item(45, 6).Value = "2021-10-05 13:41:33.908229"
$ws.Cells.Item(46, 6).Value = "2021-10-05 13:41:33.908232"
$ws.Cells.Item(47, 6).Value = "2021-10-05 13:41:33.908235"

$excel.CutCopyMode = 0
